$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'330.87"
$ws.Range("E2").Value = "'1.16%"
$ws.Range("D3").Value = "'41.08"
$ws.Range("E3").Value = "'2.33%"
$ws.Range("D4").Value = "'5.734"
$ws.Range("E4").Value = "'-2.20%"
$ws.Range("D5").Value = "'0.08130"
$ws.Range("E5").Value = "'1.43%"
$ws.Range("D6").Value = "'2.058"
$ws.Range("E6").Value = "'7.07%"
$ws.Range("D7").Value = "'8.741"
$ws.Range("E7").Value = "'0.38%"
$ws.Range("D8").Value = "'4.510"
$ws.Range("E8").Value = "'-1.53%"
$ws.Range("E9").Value = "'0.77%"
$ws.Range("D10").Value = "'0.9219"
$ws.Range("E10").Value = "'-2.11%"
$ws.Range("D11").Value = "'0.1241"
$ws.Range("E11").Value = "'-1.20%"
$ws.Range("D12").Value = "'0.1953"
$ws.Range("E12").Value = "'-0.72%"
$ws.Range("D13").Value = "'8.319"
$ws.Range("E13").Value = "'-6.25%"
$ws.Range("D14").Value = "'0.09324"
$ws.Range("E14").Value = "'2.16%"
$ws.Range("D15").Value = "'0.03665"
$ws.Range("E15").Value = "'2.91%"
$ws.Range("D16").Value = "'0.1055"
$ws.Range("E16").Value = "'9.49%"
$ws.Range("D17").Value = "'0.001302"
$ws.Range("E17").Value = "'0.35%"
$ws.Range("D18").Value = "'0.006217"
$ws.Range("E18").Value = "'0.40%"
$ws.Range("D19").Value = "'3.381"
$ws.Range("E19").Value = "'0.65%"
$ws.Range("E20").Value = "'-1.21%"
$ws.Range("D21").Value = "'0.1416"
$ws.Range("E21").Value = "'-1.11%"
$ws.Range("E22").Value = "'9.71%"
$ws.Range("D23").Value = "'0.04421"
$ws.Range("E23").Value = "'-0.67%"
$ws.Range("E24").Value = "'-0.04%"
$ws.Range("D25").Value = "'0.004356"
$ws.Range("E25").Value = "'0.95%"
$ws.Range("E26").Value = "'8.52%"
$ws.Range("D39").Value = "'0.02784"
$ws.Range("E39").Value = "'14.98%"
$ws.Range("D40").Value = "'0.05503"
$ws.Range("E40").Value = "'4.29%"
$ws.Range("D41").Value = "'0.007603"
$ws.Range("E41").Value = "'1.97%"
$ws.Range("D42").Value = "'0.009934"
$ws.Range("E42").Value = "'14.30%"
$ws.Range("D43").Value = "'0.1425"
$ws.Range("E43").Value = "'0.60%"
$ws.Range("D44").Value = "'0.002120"
$ws.Range("E44").Value = "'-0.32%"
$ws.Range("D45").Value = "'0.01185"
$ws.Range("E45").Value = "'11.54%"
$ws.Range("D46").Value = "'0.00006738"
$ws.Range("E46").Value = "'-1.69%"
$ws.Range("D47").Value = "'0.00000000750"
$ws.Range("E47").Value = "'-0.38%"
$ws.Range("D48").Value = "'0.002942"
$ws.Range("E48").Value = "'-6.73%"
$ws.Range("D49").Value = "'0.002278"
$ws.Range("E49").Value = "'59.99%"
$ws.Range("D50").Value = "'0.00002099"
$ws.Range("E50").Value = "'-0.38%"
$ws.Range("D51").Value = "'0.0001999"
$ws.Range("E51").Value = "'-0.38%"
